$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Lists (continued)"
